$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -18.84761551611473
$ws.Range("B3").Value = -1.557017206064529
$ws.Range("B4").Value = 0.3586898072596796
$ws.Range("B5").Value = 2.42130931526367
$ws.Range("B6").Value = 4.489374712105018
$ws.Range("B7").Value = 15.49025017839676
$ws.Range("B8").Value = 18.17790878329038
$ws.Range("B9").Value = 10.3675845361156
$ws.Range("B10").Value = -19.69552273777981
$ws.Range("B11").Value = -52.26237632495349
$ws.Range("B12").Value = -63.82773632146177
$ws.Range("B13").Value = -67.07750751751269
$ws.Range("B14").Value = -63.38963396911004
$ws.Range("B15").Value = -59.41941412923393
$ws.Range("B16").Value = -57.86153978682387
$ws.Range("B17").Value = -56.43112285542816
$ws.Range("B18").Value = -49.36377664711154
$ws.Range("B19").Value = -46.18761868785487
$ws.Range("B20").Value = -53.71358902473526
$ws.Range("B21").Value = -72.33047157877172
$ws.Range("B22").Value = -83.93328975729355
$ws.Range("B23").Value = -92.04840548889388
$ws.Range("B24").Value = -94.68488722073516
$ws.Range("B25").Value = -76.08070572768224
